$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TIMER Clock Frequency (TIMx_CLK)
$ws.Range("B4").Value = 50000000

# Update Prescaler
$ws.Range("B6").Value = 4

# Update Time base Required
$ws.Range("B14").Value = 0.001

# Update selected cell to B7
$ws.Range("B7").Select()

$wb.Save()
